$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New comuni subentrati rows are appended below the existing table
# (rows 2-14). Formatting for each new cell is copied from the row
# directly above (row 14) so borders / quote-prefixed-text / date
# styles match the rest of the table exactly.

# --- Row 15 : SPINEA (VE) ---
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = "SPINEA"

$ws.Range("F14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F15").Value = "VE"

$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = "'027038"

$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = "I908"

$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = "'05"

$ws.Range("D15").Value = 27
$ws.Range("D14").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("G14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G15").Value = 42991.75

# --- Row 16 : DRO (TN) ---
$ws.Range("A14").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = "'022079"

$ws.Range("B14").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B16").Value = "D371"

$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = "DRO"

$ws.Range("E14").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Value = "'04"

$ws.Range("F14").Copy()
$ws.Range("F16").PasteSpecial(-4122)
$ws.Range("F16").Value = "TN"

$ws.Range("D16").Value = 22
$ws.Range("D14").Copy()
$ws.Range("D16").PasteSpecial(-4122)

$ws.Range("G14").Copy()
$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("G16").Value = 42996.666666666664

# Match the author's final selection/active cell.
$ws.Range("F32").Select() | Out-Null
